# Menu-Languages.docx (Indonesian) translation update:
# "RPC Explorer" -> "Insight Explorer"
$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("RPC Explorer", $true, $false, $false, $false, $false, $true, 1, $false, "Insight Explorer", 2)
